# Updated car model drag coeff and added low drag configs

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("Info")

# --- Aerodynamics figures on the Info sheet ---
# Lift Coefficient CL
$wsInfo.Range("C8").Value = -1.98
# Drag Coefficient CD
$wsInfo.Range("C9").Value = -1.33
# Front Aero Distribution, now computed from a formula
$wsInfo.Range("C12").Formula = "=100-56.3"
# Frontal Area
$wsInfo.Range("C13").Value = 1.15

# --- Make "Info" the active sheet/tab and update its selection & scroll ---
$wsInfo.Activate()
$wsInfo.Range("E12").Select()

# --- "Torque Curve" keeps its own selection, just loses the scroll position ---
# (selection stays at J48, which is already the case in the source file)
